$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Sources" header to "References"
$ws.Range("G1").Value = "References"

# Widen column G to fit the new, longer "References" header/links
# (47 characters of stored OOXML width == 46.1666... in COM "characters" units)
$ws.Columns("G").ColumnWidth = 46.166666666666664

# Row 7 previously lacked the yellow highlight used by rows 3-6; bring it
# into line with the rest of the "Created Maven Project" block.
$ws.Range("A7:G7").Interior.Color = 65535

# Add the new reference link for row 7 (Dynamic Programming Implementation)
$ws.Range("G7").Value = "https://www.youtube.com/watch?v=JE0JE8ce1V0"

# Update the active selection
$ws.Range("G14").Select()
